{"js": "// Office.js (Word JavaScript API) edit script.\n// This is the body of: async (context) => { ... }\n//\n// Changes applied (per the target diff):\n//   1. Bump the generated \"Date\" paragraph's timestamp.\n//   2. Fix the \"ressources\" -> \"resources\" typo.\n//   3. Fix the \"occasionaly\" -> \"occasionally\" typo.\n//   4. Extend the file-format list (\"html, pdf, and odt\") to also\n//      mention \"docx\" (\"html, pdf, odt and docx\"), keeping the\n//      \"docx\" token styled the same way (NormalTok) as the other\n//      format tokens (html / pdf / odt).\n\nconst body = context.document.body;\n\n// 1. Update the generated-date line (the string is unique in the document).\nconst dateRanges = body.search(\"May  27, 2021 (11:54:01 PM)\", { matchCase: true });\ndateRanges.load(\"items\");\nawait context.sync();\nfor (const r of dateRanges.items) {\n  r.insertText(\"May  28, 2021 (01:53:57 AM)\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2. Fix \"ressources\" -> \"resources\" (unique occurrence).\nconst typoRanges1 = body.search(\"ressources\", { matchCase: true });\ntypoRanges1.load(\"items\");\nawait context.sync();\nfor (const r of typoRanges1.items) {\n  r.insertText(\"resources\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 3. Fix \"occasionaly\" -> \"occasionally\" (unique occurrence).\nconst typoRanges2 = body.search(\"occasionaly\", { matchCase: true });\ntypoRanges2.load(\"items\");\nawait context.sync();\nfor (const r of typoRanges2.items) {\n  r.insertText(\"occasionally\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 4. Turn \"..., and odt thanks to ...\" into \"..., odt and docx thanks to ...\".\n//    Locate the list-item paragraph describing the markdown conversion first\n//    and scope every further search to it, so the edits never touch the\n//    unrelated \"...possible, and occasionally...\" sentence, or the second,\n//    unrelated \"...highlighted thanks to Pygments.\" paragraph later on.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nlet formatsParagraph = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"markdown source code is converted to\") !== -1) {\n    formatsParagraph = p;\n    break;\n  }\n}\n\n//    a) The run reading \", and\" (right after \"pdf\") becomes just \",\".\nconst commaRanges = formatsParagraph.search(\", and\", { matchCase: true });\ncommaRanges.load(\"items\");\nawait context.sync();\nconst commaRange = commaRanges.items[0];\ncommaRange.insertText(\",\", Word.InsertLocation.replace);\nawait context.sync();\n\n//    b) Insert \"and docx \" right before \"thanks to\" (so the sentence reads\n//       \"... odt and docx thanks to ...\").\nconst thanksRanges = formatsParagraph.search(\"thanks to\", { matchCase: true });\nthanksRanges.load(\"items\");\nawait context.sync();\nconst thanksRange = thanksRanges.items[0];\nconst insertionPoint = thanksRange.getRange(Word.RangeLocation.before);\ninsertionPoint.insertText(\"and docx \", Word.InsertLocation.replace);\nawait context.sync();\n\n//    c) Give the newly-inserted \"docx\" token the same character style\n//       (\"NormalTok\") used by the other format tokens (html/pdf/odt).\nconst docxRanges = formatsParagraph.search(\"docx\", { matchCase: true });\ndocxRanges.load(\"items,text\");\nawait context.sync();\nfor (const r of docxRanges.items) {\n  if (r.text === \"docx\") {\n    r.set({ style: \"NormalTok\" });\n  }\n}\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Changes applied (per the target diff):\n#   1. Bump the generated \"Date\" paragraph's timestamp.\n#   2. Fix the \"ressources\" -> \"resources\" typo.\n#   3. Fix the \"occasionaly\" -> \"occasionally\" typo.\n#   4. Extend the file-format list (\"html, pdf, and odt\") to also\n#      mention \"docx\" (\"html, pdf, odt and docx\"), keeping the \"docx\"\n#      token styled the same way (NormalTok) as the other format\n#      tokens (html / pdf / odt).\n\n$d = $word.ActiveDocument\n\n# 1. Update the generated-date line. The string is unique in the document,\n#    so a plain whole-document Find/Replace is safe.\n$dateRange = $d.Content\n$null = $dateRange.Find.Execute(\n  \"May  27, 2021 (11:54:01 PM)\", $false, $false, $false, $false, $false,\n  $true, 1, $false, \"May  28, 2021 (01:53:57 AM)\", 2)\n\n# 2. Fix \"ressources\" -> \"resources\" (unique occurrence).\n$typoRange1 = $d.Content\n$null = $typoRange1.Find.Execute(\n  \"ressources\", $true, $false, $false, $false, $false,\n  $true, 1, $false, \"resources\", 2)\n\n# 3. Fix \"occasionaly\" -> \"occasionally\" (unique occurrence).\n$typoRange2 = $d.Content\n$null = $typoRange2.Find.Execute(\n  \"occasionaly\", $true, $false, $false, $false, $false,\n  $true, 1, $false, \"occasionally\", 2)\n\n# 4. Turn \"..., and odt thanks to ...\" into \"..., odt and docx thanks to ...\".\n#    Locate the list-item paragraph describing the markdown conversion so we\n#    never touch the unrelated \"...possible, and occasionally...\" sentence\n#    earlier in the document.\n$paras = $d.Paragraphs\n$formatsPara = $null\nfor ($i = 1; $i -le $paras.Count; $i++) {\n  $p = $paras.Item($i)\n  if ($p.Range.Text -like \"*markdown source code is converted to*\") {\n    $formatsPara = $p\n    break\n  }\n}\n\n#    a) The run reading \", and\" (right after \"pdf\") becomes just \",\": delete\n#       the \" and\" suffix only, leaving the original (unstyled) comma run\n#       untouched so \"pdf\" keeps its own NormalTok-only run.\n$commaScope = $formatsPara.Range.Duplicate\n$null = $commaScope.Find.Execute(\", and\")\n$suffix = $d.Range($commaScope.Start + 1, $commaScope.End)\n$suffix.Delete()\n\n#    b) Insert \"and docx \" right before \"thanks to\" (so the sentence reads\n#       \"... odt and docx thanks to ...\").\n$thanksScope = $formatsPara.Range.Duplicate\n$null = $thanksScope.Find.Execute(\"thanks to\")\n$insertionPoint = $d.Range($thanksScope.Start, $thanksScope.Start)\n$insertionPoint.InsertBefore(\"and docx \")\n\n#    c) Give the newly-inserted \"docx\" token the same character style\n#       (\"NormalTok\") used by the other format tokens (html/pdf/odt).\n$docxScope = $formatsPara.Range.Duplicate\n$null = $docxScope.Find.Execute(\"docx\")\n$docxScope.Style = \"NormalTok\"\n"}
